$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Num($ref, $val) {
    $ws.Range($ref).Value = $val
}

function Set-EmptyText($ref) {
    # Force a text entry of an empty string (matches the source file's
    # empty inlineStr cells: Text type, raw value ""), then strip the
    # auto-applied quote-prefix style so no stray formatting is left behind.
    $ws.Range($ref).Value = "'"
    $ws.Range($ref).Style = "Normal"
}

# --- Row 2: RM 2 --- C2 goes from missing -> 14.9
Set-Num "C2" 14.9

# --- Row 3: RM 8 --- D3 goes from -14.2 -> missing
Set-EmptyText "D3"

# --- Row 4: RM 9 --- D4 goes from missing -> -15.4
Set-Num "D4" -15.4

# --- Row 6: RM 21 --- C6 goes from 15.1 -> missing
Set-EmptyText "C6"

# --- Row 8: RM 38 --- D8 goes from -13.9 -> missing
Set-EmptyText "D8"

# --- Row 9: RM 42 --- D9 goes from -14.5 -> missing
Set-EmptyText "D9"

# --- Row 12: RM 81 --- C12 goes from missing -> 12.5
Set-Num "C12" 12.5

# --- Row 14: RM 90 --- C14 goes from 14.4 -> missing
Set-EmptyText "C14"

# --- Row 15: RM 95 --- D15 goes from missing -> -15.2
Set-Num "D15" -15.2

# --- Row 18: RM 120 --- D18 goes from missing -> -15.2
Set-Num "D18" -15.2

# --- Row 19: RM 125 --- D19 goes from -15.5 -> missing
Set-EmptyText "D19"

# --- Row 20: RM 134 --- C20 goes from missing -> 12.5
Set-Num "C20" 12.5

# --- Row 21: RM 135 --- C21 goes from missing -> 12.7
Set-Num "C21" 12.7

# --- Row 22: RM 138 --- D22 goes from -15.4 -> missing
Set-EmptyText "D22"

# --- Row 23: RM 140 --- C23 goes from 12.2 -> missing; D23 goes from missing -> -13.9
Set-EmptyText "C23"
Set-Num "D23" -13.9

# --- Row 24: RM 142a --- C24 goes from 12.7 -> missing
Set-EmptyText "C24"

# --- Row 25: RM 145 --- D25 goes from missing -> -15.5
Set-Num "D25" -15.5

# Two rows are dropped from the bottom block: "RM 232" (row 26) and
# "SC 92" (row 28). Deleting the higher-numbered row first keeps the
# row-26 reference valid for the second delete.
$ws.Range("A28").EntireRow.Delete()
$ws.Range("A26").EntireRow.Delete()

# After the deletes, rows 26-33 hold (in order): SC 5, SC 101, SC 105,
# SC 119, SC 120, SC 132, SC 193, SC 232 - but several of their values
# differ from what simply shifted up, so re-assert the final values.

# Row 26: SC 5
Set-Num "B26" -20.2
Set-Num "C26" 10.8
Set-Num "D26" -13.8
Set-Num "E26" -5
Set-Num "F26" 17.38

# Row 27: SC 101
Set-EmptyText "B27"
Set-Num "C27" 10
Set-EmptyText "D27"
Set-Num "E27" -10
Set-Num "F27" 17

# Row 28: SC 105
Set-EmptyText "B28"
Set-Num "C28" 11.1
Set-Num "D28" -13.7
Set-Num "E28" -5.9
Set-Num "F28" 17.44

# Row 29: SC 119
Set-Num "B29" -19.5
Set-Num "C29" 11.2
Set-EmptyText "D29"
Set-Num "E29" -6.8
Set-Num "F29" 18.06

# Row 30: SC 120
Set-Num "B30" -19.7
Set-Num "C30" 11.4
Set-Num "D30" -13.6
Set-Num "E30" -5.7
Set-Num "F30" 16.89

# Row 31: SC 132
Set-EmptyText "B31"
Set-Num "C31" 15.3
Set-Num "D31" -13.7
Set-Num "E31" -8.1
Set-Num "F31" 17.18

# Row 32: SC 193
Set-EmptyText "B32"
Set-Num "C32" 10.5
Set-Num "D32" -14.7
Set-Num "E32" -6.4
Set-Num "F32" 17.39

# Row 33: SC 232
Set-Num "B33" -19.5
Set-Num "C33" 10.4
Set-Num "D33" -14.1
Set-Num "E33" -10.7
Set-Num "F33" 17.53
